$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the missing "F" (Friday) hours for week 1
$ws.Range("G2").Value = 4.5

# Update the active cell selection to match the latest saved state
$ws.Range("H13").Select()

$wb.Save()
